$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E5").Value = 12.66959999999999
$ws.Range("E6").Value = 12.5027
$ws.Range("D7").Value = -7.335299999999994
$ws.Range("B10").Value = 8.653300000000005
$ws.Range("B12").Value = 6.194400000000003
$ws.Range("C13").Value = -13.20469999999999
$ws.Range("B18").Value = 4.781500000000007
$ws.Range("D20").Value = -8.255500000000001
